$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "krthik234@gmail.com"
$ws.Range("B3").Value = "password"
$ws.Outline.ShowLevels(2, 1)
$ws.Range("C10").Select()
